# "revert to version 12, small code cleanup"
# Adds a new evaluation row (commit a523aff results + new commit fa6507f
# entry about switching the augmentation pipeline to albumentations) to the
# "Тесты" log sheet, and moves the sheet's scroll/selection forward.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Тесты")
$ws.Activate()

# --- Row 25 (test #24, commit a523aff): fill in the results that were
#     still pending in the previous version of the sheet. ---
$ws.Range("H25").Value = "Train IoU: 0.50, Val IoU: 0.45. "
$ws.Range("I25").Value = "a523aff"

# --- Row 26 (test #25, commit fa6507f): brand new entry. ---
$ws.Range("B26").Value = 1
$ws.Range("C26").Value = 40
$ws.Range("D26").Value = 17
$ws.Range("F26").Value = "Изменена система аугментаций"
$ws.Range("G26").Value = "Вместо старой системы аугментации используется улучшенный набор аугментаций из albumentations"
$ws.Range("H26").Value = "Train IoU: 0.43, Val IoU: 0.46. Нуждается в доработке"
$ws.Range("I26").Value = "fa6507f"

# The new row wraps onto two lines in the "Параметры предобработки..." /
# "Результаты" columns, so give it the taller row height Excel would have
# picked via auto-fit.
$ws.Rows.Item(26).RowHeight = 30

# Scroll the view down / select the new row's "Результаты" cell, same as
# the author left the sheet.
$ws.Range("A16").Select()
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("H27").Select()
